$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.131.29"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "1.551.64"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "287.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3806"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3312"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07369"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.806"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.708"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.07%  "

$ws.Range("D16").Value = "1.567.76"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001070"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06650"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.53%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.373"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.43%  "

$ws.Range("D24").Value = "22.135.68"
$ws.Range("E24").Value = "  -1.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.291"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.532"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.931"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("B30").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C30").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.764.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.01%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.089"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.869"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.910"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08191"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.301"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06280"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02314"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.68%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.277"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2142"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.224"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.29%  "

$ws.Range("B43").Value = "Frax"
$ws.Range("C43").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6039"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.736"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5840"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.958"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.171"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07026"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.26%  "
